$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row containing "Previous" -> "Previously", height 607 -> 624 twips (30.35 -> 31.2 pt)
$rowPrev = $t.Rows.Item(31)
$rowPrev.Height = 31.2
$rowPrev.Cells.Item(1).Range.Find.Execute("Previous", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Previously", 2)

# Row containing "Current" -> "Currently", height 607 -> 624 twips (30.35 -> 31.2 pt)
$rowCurr = $t.Rows.Item(32)
$rowCurr.Height = 31.2
$rowCurr.Cells.Item(1).Range.Find.Execute("Current", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Currently", 2)
